# Lab 3 rubric grade change + minor fixes
#
# Semantic summary of the authored edit (reconstructed from the OOXML diff):
#   - The "Possible" points in column D (rows 6,8-18) were re-totaled from
#     50 down to 40 on both the "Rubric" and "Score" sheets. The "Score"
#     column E on the "Score" sheet mirrors the same new values.
#   - The "Comment" column (G) on the "Score" sheet had Wrap Text turned on
#     for the score rows (G6:G18).
#   - The D20/E20 "Total" cells are =SUM(...) formulas and recompute on
#     their own once the inputs above change.

$wb = $excel.ActiveWorkbook

# New "Possible" point values for rows 6, 8..18 (row 7 has no score - it's
# the "Part 2" sub-header). Old -> New: 10->8, 4->3, 5->4, 4->3, 6->5, 6->4, 6->4
$newValues = @{
    6  = 8
    8  = 2
    9  = 3
    10 = 1
    11 = 4
    12 = 3
    13 = 2
    14 = 5
    15 = 4
    16 = 2
    17 = 4
    18 = 2
}

# --- "Rubric" sheet: column D only ---
# The now-unused "Score" (E) column and the leftover blank comment (G6)
# cell next to the possible-points range are cleared out entirely (they
# carried only stale formatting, no data).
$wsRubric = $wb.Worksheets.Item("Rubric")
$wsRubric.Range("E6:E18").Clear()
$wsRubric.Range("G6").Clear()
foreach ($row in $newValues.Keys) {
    $cell = $wsRubric.Cells.Item($row, 4)
    $cell.ClearFormats()
    $cell.Value = $newValues[$row]
}

# --- "Score" sheet: columns D and E mirror each other ---
$wsScore = $wb.Worksheets.Item("Score")
foreach ($row in $newValues.Keys) {
    $dCell = $wsScore.Cells.Item($row, 4)
    $dCell.ClearFormats()
    $dCell.Value = $newValues[$row]

    $eCell = $wsScore.Cells.Item($row, 5)
    $eCell.ClearFormats()
    $eCell.Value = $newValues[$row]
}

# Minor fix: wrap text was turned on for the "Comment" column next to the
# score rows on the "Score" sheet.
$wsScore.Range("G6:G18").WrapText = $true

# Leave the selection on the range that was just edited, matching where the
# cursor was left after the grade change.
$wsRubric.Activate() | Out-Null
$wsRubric.Range("D6:D18").Select() | Out-Null

$wsScore.Activate() | Out-Null
$wsScore.Range("H14").Select() | Out-Null

$excel.Calculate()
